$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (price) cells to text so that dotted/decimal-looking strings
# are preserved exactly as text, matching the original inline-string cell type.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "60.263.84"
$ws.Range("E2").Value = "  -5.97%  "
$ws.Range("D3").Value = "3.254.54"
$ws.Range("E3").Value = "  -6.31%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "559.51"
$ws.Range("E5").Value = "  -4.17%  "
$ws.Range("D6").Value = "126.36"
$ws.Range("E6").Value = "  -3.62%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "3.252.80"
$ws.Range("E8").Value = "  -6.25%  "
$ws.Range("D9").Value = "0.472"
$ws.Range("E9").Value = "  -2.01%  "
$ws.Range("D10").Value = "7.25"
$ws.Range("E10").Value = "  -4.47%  "
$ws.Range("E11").Value = "  -4.86%  "
$ws.Range("D12").Value = "0.369"
$ws.Range("E12").Value = "  -4.19%  "
$ws.Range("D13").Value = "3.823.24"
$ws.Range("E13").Value = "  -6.04%  "
$ws.Range("E14").Value = "  -0.06%  "
$ws.Range("D15").Value = "3.268.49"
$ws.Range("E15").Value = "  -6.06%  "
$ws.Range("E16").Value = "  -6.12%  "
$ws.Range("D17").Value = "60.434.94"
$ws.Range("E17").Value = "  -5.71%  "
$ws.Range("D18").Value = "24.14"
$ws.Range("E18").Value = "  -0.50%  "
$ws.Range("D19").Value = "5.59"
$ws.Range("E19").Value = "  -1.38%  "
$ws.Range("D20").Value = "13.16"
$ws.Range("E20").Value = "  -1.63%  "
$ws.Range("D21").Value = "8.94"
$ws.Range("E21").Value = "  -9.98%  "
$ws.Range("D22").Value = "349.84"
$ws.Range("E22").Value = "  -9.18%  "
$ws.Range("E23").Value = "  -3.09%  "
$ws.Range("D24").Value = "0.999"
$ws.Range("E24").Value = "  -0.12%  "
$ws.Range("D25").Value = "3.393.52"
$ws.Range("E25").Value = "  -6.17%  "
$ws.Range("D26").Value = "69.07"
$ws.Range("E26").Value = "  -7.53%  "
$ws.Range("E27").Value = "  -4.82%  "
$ws.Range("E28").Value = "  +0.29%  "
$ws.Range("D29").Value = "7.19"
$ws.Range("E29").Value = "  +1.30%  "
$ws.Range("E30").Value = "  -1.58%  "
$ws.Range("D31").Value = "7.76"
$ws.Range("E31").Value = "  -1.73%  "
$ws.Range("D32").Value = "2.08"
$ws.Range("E32").Value = "  -6.39%  "
$ws.Range("E33").Value = "  +0.03%  "
$ws.Range("D34").Value = "0.148"
$ws.Range("E34").Value = "  -2.08%  "
$ws.Range("D35").Value = "3.291.46"
$ws.Range("E35").Value = "  -6.12%  "
$ws.Range("D36").Value = "22.48"
$ws.Range("E36").Value = "  -1.89%  "
$ws.Range("D37").Value = "5.20"
$ws.Range("E37").Value = "  +0.75%  "
$ws.Range("D38").Value = "6.74"
$ws.Range("E38").Value = "  +0.17%  "
$ws.Range("E39").Value = "  -1.86%  "
$ws.Range("D40").Value = "157.92"
$ws.Range("E40").Value = "  -2.85%  "
$ws.Range("D41").Value = "0.0747"
$ws.Range("E41").Value = "  -3.26%  "
$ws.Range("E42").Value = "  +0.17%  "
$ws.Range("D43").Value = "40.79"
$ws.Range("E43").Value = "  -1.34%  "
$ws.Range("D44").Value = "4.30"
$ws.Range("E44").Value = "  +0.56%  "
$ws.Range("D45").Value = "0.732"
$ws.Range("E45").Value = "  -7.88%  "
$ws.Range("E47").Value = "  -4.58%  "
$ws.Range("D48").Value = "22.43"
$ws.Range("E48").Value = "  -3.77%  "
$ws.Range("D49").Value = "6.62"
$ws.Range("E49").Value = "  -0.74%  "
$ws.Range("D50").Value = "0.854"
$ws.Range("E50").Value = "  -4.75%  "
$ws.Range("D51").Value = "21.13"
$ws.Range("E51").Value = "  +3.82%  "

# Restore default styling on column D so no extra formatting is introduced.
$ws.Range("D2:D51").Style = "Normal"
